$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-10: Sending cluster / Ligand / Receptor / Target cluster pairs
# for the three clusters (ECs, FAPs, sCs) crossed with the Epo-Ephb4 pair.
$rows = @(
    @("ECs", "Epo", "Ephb4", "ECs", 2, 0.6666666666666666, 0.188903, 0.566709, 0.3560183867666543, 0.3560183867666543, 2, 0.6666666666666666, 23.746319, 71.238957, 0.7135031414879517, 0.7135031414879517, 4.485750898057001, 40.371758082513, 0.2540202373854805, 0.2540202373854804),
    @("ECs", "Epo", "Ephb4", "FAPs", 2, 0.6666666666666666, 0.188903, 0.566709, 0.3560183867666543, 0.3560183867666543, 3, 1, 4.865208333333334, 14.595625, 0.1461844014571983, 0.1461844014571983, 0.9190524497916669, 8.271472048125002, 0.05204433477724069, 0.05204433477724067),
    @("ECs", "Epo", "Ephb4", "sCs", 2, 0.6666666666666666, 0.188903, 0.566709, 0.3560183867666543, 0.3560183867666543, 3, 1, 4.669782333333333, 14.009347, 0.1403124570548501, 0.1403124570548501, 0.8821358921136667, 7.939223029023, 0.04995381460393319, 0.04995381460393318),
    @("FAPs", "Epo", "Ephb4", "ECs", 3, 1, 0.2247596666666667, 0.6742790000000001, 0.4235960992513493, 0.4235960992513493, 2, 0.6666666666666666, 23.746319, 71.238957, 0.7135031414879517, 0.7135031414879517, 5.337214743000334, 48.03493268700301, 0.3022371475378799, 0.3022371475378799),
    @("FAPs", "Epo", "Ephb4", "FAPs", 3, 1, 0.2247596666666667, 0.6742790000000001, 0.4235960992513493, 0.4235960992513493, 3, 1, 4.865208333333334, 14.595625, 0.1461844014571983, 0.1461844014571983, 1.093502603263889, 9.841523429375002, 0.06192314222866246, 0.06192314222866244),
    @("FAPs", "Epo", "Ephb4", "sCs", 3, 1, 0.2247596666666667, 0.6742790000000001, 0.4235960992513493, 0.4235960992513493, 3, 1, 4.669782333333333, 14.009347, 0.1403124570548501, 0.1403124570548501, 1.049578720645889, 9.446208485813, 0.05943580948480695, 0.05943580948480695),
    @("sCs", "Epo", "Ephb4", "ECs", 2, 0.6666666666666666, 0.1169363333333333, 0.350809, 0.2203855139819964, 0.2203855139819964, 2, 0.6666666666666666, 23.746319, 71.238957, 0.7135031414879517, 0.7135031414879517, 2.776807474023666, 24.991267266213, 0.1572457565645913, 0.1572457565645913),
    @("sCs", "Epo", "Ephb4", "FAPs", 2, 0.6666666666666666, 0.1169363333333333, 0.350809, 0.2203855139819964, 0.2203855139819964, 3, 1, 4.865208333333334, 14.595625, 0.1461844014571983, 0.1461844014571983, 0.5689196234027778, 5.120276610625, 0.03221692445129515, 0.03221692445129515),
    @("sCs", "Epo", "Ephb4", "sCs", 2, 0.6666666666666666, 0.1169363333333333, 0.350809, 0.2203855139819964, 0.2203855139819964, 3, 1, 4.669782333333333, 14.009347, 0.1403124570548501, 0.1403124570548501, 0.5460672235247777, 4.914605011722999, 0.03092283296610993, 0.03092283296610993)
)

$nRows = $rows.Count
$nCols = 20
$data = New-Object 'object[,]' $nRows, $nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A2:T10").Value = $data
